# Appends 35 new log rows (Id 442-476) to the "logs" worksheet, continuing
# directly after the existing last row (Id 441 / row 284), matching the
# date-formatted Column D style already used by the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new rows: Id, Description(col B), Message(col C), Date(col D, serial)
$newRows = @(
    @(442, "Info", "Tüm personeller listelendi", 45597),
    @(443, "Info", "Tüm personeller listelendi", 45597),
    @(444, "Info", "Kiralamalar listelendi", 45597),
    @(445, "Info", "Tüm departmanlar listelendi", 45597),
    @(446, "Info", "Tüm departmanlar listelendi", 45597),
    @(447, "Info", "Tüm personeller listelendi", 45597),
    @(448, "Info", "Kiralamalar listelendi", 45597),
    @(449, "Info", "Tüm departmanlar listelendi", 45597),
    @(450, "Info", "Tüm personeller listelendi", 45597),
    @(451, "Info", "Kiralamalar listelendi", 45597),
    @(452, "Info", "Tüm personeller listelendi", 45597),
    @(453, "Info", "Tüm personeller listelendi", 45597),
    @(454, "Info", "Tüm personeller listelendi", 45597),
    @(455, "Info", "Tüm personeller listelendi", 45597),
    @(456, "Info", "Kiralamalar listelendi", 45597),
    @(457, "Info", "Tüm personeller listelendi", 45597),
    @(458, "Info", "Tüm personeller listelendi", 45597),
    @(459, "Info", "Kiralamalar listelendi", 45597),
    @(460, "Info", "Tüm personeller listelendi", 45599),
    @(461, "Info", "Kiralamalar listelendi", 45599),
    @(462, "Info", "Kiralamalar listelendi", 45599),
    @(463, "Info", "Tüm departmanlar listelendi", 45599),
    @(464, "Info", "Tüm departmanlar listelendi", 45599),
    @(465, "Info", "Tüm personeller listelendi", 45599),
    @(466, "Info", "Tüm personeller listelendi", 45599),
    @(467, "Info", "Tüm personeller listelendi", 45599),
    @(468, "Info", "Tüm personeller listelendi", 45599),
    @(469, "Info", "Kiralamalar listelendi", 45600),
    @(470, "Info", "Kiralamalar listelendi", 45600),
    @(471, "Info", "İd değerine göre kiralama listelendi", 45600),
    @(472, "Info", "İd değerine göre kiralama listelendi", 45600),
    @(473, "Error", "İd değerine gmre kiralama listelenmesinde hata oluştu", 45600),
    @(474, "Info", "İd değerine göre kiralama listelendi", 45600),
    @(475, "Info", "Kiralamalar listelendi", 45600),
    @(476, "Info", "Tüm departmanlar listelendi", 45600)
)

# Existing data ends at row 284 (last Id = 441); new rows continue right after it.
$lastRow = $ws.UsedRange.Rows.Count
$dateFormat = $ws.Cells.Item($lastRow, 4).NumberFormat

$r = $lastRow + 1
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 4).NumberFormat = $dateFormat
    $r = $r + 1
}
